# Leave Card update - 12/22/2023 10:59 AM
# Adds the July-2023..Mar-2024 VL/SL earning entries, a Dec-2023 VL usage
# entry, a "2024" year-divider row, the GSIS "entrance to duty" date, and
# grows the leave table by one row so the special bottom border stays on
# the true last row.

$wb = $excel.ActiveWorkbook
$leave = $wb.Worksheets.Item("LEAVE CREDITS")
$conv  = $wb.Worksheets.Item("CONVERTION")

# ---------------------------------------------------------------------
# LEAVE CREDITS sheet
# ---------------------------------------------------------------------

# ENTRANCE TO DUTY date (merged F3:G3) = 07/03/2023
$leave.Range("F3").Value = 45110

# Monthly VL/SL earnings, July 2023 - Nov 2023 (1.167 for the partial
# first month, 1.25/month afterwards)
$leave.Range("A11").Value = 45110   # 07/03/2023
$leave.Range("C11").Value = 1.167

$leave.Range("A12").Value = 45169   # 08/31/2023
$leave.Range("C12").Value = 1.25

$leave.Range("A13").Value = 45199   # 09/30/2023
$leave.Range("C13").Value = 1.25

$leave.Range("A14").Value = 45230   # 10/31/2023
$leave.Range("C14").Value = 1.25

$leave.Range("A15").Value = 45260   # 11/30/2023
$leave.Range("C15").Value = 1.25

# Row 17 becomes the "2024" year-divider label: force text (not a date
# serial) and pick up the bold/banded divider style used by row 10.
$leave.Range("A17").Value = "'2024"
$leave.Range("A10").Copy()
$leave.Range("A17").PasteSpecial(-4122)

# December 2023 row also records a 6-day VL usage (absence w/ pay)
$leave.Range("A16").Value = 45291   # 12/31/2023
$leave.Range("B16").Value = "VL(6-0-0)"
$leave.Range("D16").Value = 6
$leave.Range("K16").Value = "12/12-17/2023"

# Continue the monthly earnings into 2024
$leave.Range("A18").Value = 45322   # 01/31/2024
$leave.Range("A19").Value = 45351   # 02/29/2024
$leave.Range("A20").Value = 45382   # 03/31/2024

# The divider row moved from 10/23/36/49/62/75 to 10/17/24/37/50/63/76,
# so the alternating band style on column A needs to follow it down by
# one slot for each subsequent divider.
$leave.Range("A11").Copy()
$leave.Range("A23").PasteSpecial(-4122)
$leave.Range("A36").PasteSpecial(-4122)
$leave.Range("A49").PasteSpecial(-4122)
$leave.Range("A62").PasteSpecial(-4122)
$leave.Range("A75").PasteSpecial(-4122)

$leave.Range("A10").Copy()
$leave.Range("A24").PasteSpecial(-4122)
$leave.Range("A37").PasteSpecial(-4122)
$leave.Range("A50").PasteSpecial(-4122)
$leave.Range("A63").PasteSpecial(-4122)
$leave.Range("A76").PasteSpecial(-4122)

# Grow Table15 by one row so the heavier bottom border stays on the
# genuine last row of the table (old row 134 -> new row 135; the new
# row 134 is a plain interior row cloned from row 133).
$tbl = $leave.ListObjects.Item("Table15")
$tbl.ListRows.Add() | Out-Null

$leave.Range("A134:K134").Copy()
$leave.Range("A135:K135").PasteSpecial(-4122)
$leave.Range("G135").Formula = "=IF(ISBLANK(Table15[[#This Row],[EARNED]])," + """" + """" + ",Table15[[#This Row],[EARNED]])"

$leave.Range("A133:K133").Copy()
$leave.Range("A134:K134").PasteSpecial(-4122)
$leave.Range("G134").Formula = "=IF(ISBLANK(Table15[[#This Row],[EARNED]])," + """" + """" + ",Table15[[#This Row],[EARNED]])"

$leave.Range("K16").Select()

# ---------------------------------------------------------------------
# CONVERTION sheet - daily-earn calculator: entered start day "3"
# ---------------------------------------------------------------------
$conv.Range("J3").Value = 3
$conv.Range("L3").Select()
